$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "中国卫星"
$ws.Range("B2").Value = "中国卫星"
$ws.Range("C2").Value = "航天发展"
$ws.Range("A3").Value = "锋龙股份"
$ws.Range("B3").Value = "锋龙股份"
$ws.Range("C3").Value = "锋龙股份"
$ws.Range("A4").Value = "航天发展"
$ws.Range("B4").Value = "航天发展"
$ws.Range("C4").Value = "中国卫星"
$ws.Range("A5").Value = "神剑股份"
$ws.Range("B5").Value = "航天电子"
$ws.Range("C5").Value = "天际股份"
$ws.Range("A6").Value = "再升科技"
$ws.Range("B6").Value = "神剑股份"
$ws.Range("C6").Value = "神剑股份"
$ws.Range("A7").Value = "航天电子"
$ws.Range("B7").Value = "东百集团"
$ws.Range("C7").Value = "东百集团"
$ws.Range("A8").Value = "东百集团"
$ws.Range("B8").Value = "海南发展"
$ws.Range("C8").Value = "西部材料"
$ws.Range("A9").Value = "超捷股份"
$ws.Range("B9").Value = "再升科技"
$ws.Range("C9").Value = "平潭发展"
$ws.Range("A10").Value = "天际股份"
$ws.Range("B10").Value = "昊志机电"
$ws.Range("C10").Value = "航天电子"
$ws.Range("A11").Value = "平潭发展"
$ws.Range("B11").Value = "航天动力"
$ws.Range("C11").Value = "再升科技"
$ws.Range("A12").Value = "西部材料"
$ws.Range("B12").Value = "西部材料"
$ws.Range("C12").Value = "百大集团"
$ws.Range("A13").Value = "昊志机电"
$ws.Range("B13").Value = "平潭发展"
$ws.Range("C13").Value = "金风科技"
$ws.Range("A14").Value = "华菱线缆"
$ws.Range("B14").Value = "天际股份"
$ws.Range("C14").Value = "通宇通讯"
$ws.Range("A15").Value = "海南发展"
$ws.Range("B15").Value = "通宇通讯"
$ws.Range("C15").Value = "华菱线缆"
$ws.Range("A16").Value = "通宇通讯"
$ws.Range("B16").Value = "三花智控"
$ws.Range("C16").Value = "永辉超市"
$ws.Range("A17").Value = "金风科技"
$ws.Range("B17").Value = "超捷股份"
$ws.Range("C17").Value = "国晟科技"
$ws.Range("A18").Value = "百大集团"
$ws.Range("B18").Value = "金风科技"
$ws.Range("C18").Value = "海南发展"
$ws.Range("A19").Value = "顺灏股份"
$ws.Range("B19").Value = "九鼎新材"
$ws.Range("C19").Value = "国风新材"
$ws.Range("A20").Value = "三花智控"
$ws.Range("B20").Value = "华菱线缆"
$ws.Range("C20").Value = "龙洲股份"
$ws.Range("A21").Value = "创元科技"
$ws.Range("B21").Value = "广联航空"
$ws.Range("C21").Value = "超捷股份"
